$d = $word.ActiveDocument

$replacements = @(
    @("165×8=", "142×3="),
    @("331×6=", "866×9="),
    @("823×5=", "603×5="),
    @("913×6=", "698×2="),
    @("971×8=", "881×2="),
    @("834×3=", "133×5="),
    @("603×8=", "129×2="),
    @("583×9=", "146×7="),
    @("622×6=", "722×2="),
    @("649×4=", "587×6="),
    @("586×2=", "953×4="),
    @("789×5=", "155×9="),
    @("477×9=", "987×9="),
    @("453×8=", "954×3="),
    @("677×8=", "145×2="),
    @("708×9=", "550×4="),
    @("236×3=", "708×9="),
    @("657×8=", "624×4="),
    @("678×3=", "350×2="),
    @("905×4=", "388×5="),
    @("816×2=", "778×4="),
    @("468×7=", "986×3="),
    @("587×9=", "465×6="),
    @("377×3=", "401×2="),
    @("307×9=", "908×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
